$d = $word.ActiveDocument

# Locate the run of text that needs to be split into the new
# "<comment>c_079v_02</comment>" rendition-spec runs followed by the
# remaining "Contrepoison contre la " text.
$find = $d.Content
$found = $find.Find.Execute(" Contrepoison contre la ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target text ' Contrepoison contre la '"
}

# Re-seat a plain Range at the same offsets: InsertXML must be called on a
# freshly-bound Range (the Range returned in-place by Find gets confused
# about its own bounds once content is spliced into the document).
$targetStart = $find.Start
$targetEnd = $find.End
$r = $d.Range($targetStart, $targetEnd)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:color w:val="000000"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> &lt;comment</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">&gt;</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">c_079v_02&lt;/comment&gt;</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Contrepoison contre la </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
